$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") rows 2..205: update date serial value 45178 -> 45179
$ws.Range("C2:C205").Value = 45179
